# Scheduled market-data refresh: updates currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) across the per-class leve-profit sheets with
# freshly pulled Universalis price data. Values only - no formulas/formatting.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 18
$ws.Range("H18").Value = 1190
$ws.Range("I18").Value = 1312.5
$ws.Range("K18").Value = 1312.5
$ws.Range("M18").Value = -1028.5

# row 40
$ws.Range("H40").Value = 1798.5
$ws.Range("I40").Value = 1415.3334
$ws.Range("J40").Value = 1962.7142
$ws.Range("K40").Value = 1415.3334
$ws.Range("L40").Value = 1962.7142
$ws.Range("M40").Value = -1240.3334
$ws.Range("N40").Value = -2312.7142

# row 74
$ws.Range("H74").Value = 44924.25
$ws.Range("I74").Value = 3432.3333
$ws.Range("J74").Value = 169400
$ws.Range("K74").Value = 3432.3333
$ws.Range("L74").Value = 169400
$ws.Range("M74").Value = -2496.3333
$ws.Range("N74").Value = -171272

# row 77
$ws.Range("H77").Value = 44924.25
$ws.Range("I77").Value = 3432.3333
$ws.Range("J77").Value = 169400
$ws.Range("K77").Value = 17161.6665
$ws.Range("L77").Value = 847000
$ws.Range("M77").Value = -12481.6665
$ws.Range("N77").Value = -856360

# row 80
$ws.Range("H80").Value = 1008.44446
$ws.Range("J80").Value = 2495.5
$ws.Range("L80").Value = 7486.5
$ws.Range("N80").Value = -9482.5

# row 83
$ws.Range("H83").Value = 1008.44446
$ws.Range("J83").Value = 2495.5
$ws.Range("L83").Value = 22459.5
$ws.Range("N83").Value = -32443.5

# row 98
$ws.Range("H98").Value = 1247.3334
$ws.Range("I98").Value = 696.7
$ws.Range("J98").Value = 4000.5
$ws.Range("K98").Value = 696.7
$ws.Range("L98").Value = 4000.5
$ws.Range("M98").Value = 801.3
$ws.Range("N98").Value = -6996.5

# row 111
$ws.Range("H111").Value = 449
$ws.Range("I111").Value = 449
$ws.Range("K111").Value = 1347
$ws.Range("M111").Value = 1720

# row 122
$ws.Range("H122").Value = 1247.3334
$ws.Range("I122").Value = 696.7
$ws.Range("J122").Value = 4000.5
$ws.Range("K122").Value = 2090.1
$ws.Range("L122").Value = 12001.5
$ws.Range("M122").Value = 359.8999999999996
$ws.Range("N122").Value = -16901.5

# row 137
$ws.Range("H137").Value = 1927.1666
$ws.Range("I137").Value = 1644.2632
$ws.Range("J137").Value = 3002.2
$ws.Range("K137").Value = 4932.7896
$ws.Range("L137").Value = 9006.599999999999
$ws.Range("M137").Value = -2382.7896
$ws.Range("N137").Value = -14106.6

$ws = $wb.Worksheets.Item("ARM")
# row 5
$ws.Range("H5").Value = 100.36364
$ws.Range("I5").Value = 95.40000000000001
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 95.40000000000001
$ws.Range("L5").Value = 150
$ws.Range("M5").Value = 16.59999999999999
$ws.Range("N5").Value = -374

# row 13
$ws.Range("H13").Value = 17873.25
$ws.Range("J13").Value = 17873.25
$ws.Range("L13").Value = 17873.25
$ws.Range("N13").Value = -18161.25

# row 61
$ws.Range("H61").Value = 2439.2
$ws.Range("I61").Value = 2799.25
$ws.Range("K61").Value = 2799.25
$ws.Range("M61").Value = -2587.25

# row 122
$ws.Range("H122").Value = 1766.3334
$ws.Range("I122").Value = 1199.5
$ws.Range("K122").Value = 3598.5
$ws.Range("M122").Value = -1148.5

# row 136
$ws.Range("H136").Value = 2439.2
$ws.Range("I136").Value = 2799.25
$ws.Range("K136").Value = 8397.75
$ws.Range("M136").Value = -5847.75

$ws = $wb.Worksheets.Item("BSM")
# row 4
$ws.Range("H4").Value = 100.36364
$ws.Range("I4").Value = 95.40000000000001
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 95.40000000000001
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = 19.59999999999999
$ws.Range("N4").Value = -380

# row 26
$ws.Range("H26").Value = 16000
$ws.Range("I26").Value = 16000
$ws.Range("K26").Value = 16000
$ws.Range("M26").Value = -15708

# row 134
$ws.Range("H134").Value = 13998
$ws.Range("I134").Value = 13998
$ws.Range("K134").Value = 41994
$ws.Range("M134").Value = -39459

$ws = $wb.Worksheets.Item("CRP")
# row 12
$ws.Range("H12").Value = 3808.5
$ws.Range("I12").Value = 282.33334
$ws.Range("J12").Value = 7334.6665
$ws.Range("K12").Value = 282.33334
$ws.Range("L12").Value = 7334.6665
$ws.Range("M12").Value = -112.33334
$ws.Range("N12").Value = -7674.6665

# row 60
$ws.Range("H60").Value = 25000
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# row 23
$ws.Range("H23").Value = 451.5
$ws.Range("I23").Value = 487.25
$ws.Range("K23").Value = 1461.75
$ws.Range("M23").Value = -1226.75

# row 55
$ws.Range("H55").Value = 1416.6364
$ws.Range("I55").Value = 495
$ws.Range("J55").Value = 1508.8
$ws.Range("K55").Value = 1485
$ws.Range("L55").Value = 4526.4
$ws.Range("M55").Value = -1308
$ws.Range("N55").Value = -4880.4

# row 113
$ws.Range("H113").Value = 837.52
$ws.Range("I113").Value = 651.3
$ws.Range("J113").Value = 961.6667
$ws.Range("K113").Value = 1953.9
$ws.Range("L113").Value = 2885.0001
$ws.Range("M113").Value = 216.1000000000001
$ws.Range("N113").Value = -7225.0001

$ws = $wb.Worksheets.Item("GSM")
# row 99
$ws.Range("H99").Value = 7468.25
$ws.Range("I99").Value = 7468.25
$ws.Range("K99").Value = 7468.25
$ws.Range("M99").Value = -5222.25

# row 122
$ws.Range("H122").Value = 3323.5
$ws.Range("I122").Value = 3199.1
$ws.Range("K122").Value = 9597.299999999999
$ws.Range("M122").Value = -7147.299999999999

# row 126
$ws.Range("H126").Value = 6374.75
$ws.Range("I126").Value = 8133
$ws.Range("K126").Value = 24399
$ws.Range("M126").Value = -21929

$ws = $wb.Worksheets.Item("LTW")
# row 17
$ws.Range("H17").Value = 25000
$ws.Range("J17").Value = 25000
$ws.Range("L17").Value = 25000
$ws.Range("N17").Value = -25340

# row 19
$ws.Range("H19").Value = 500425
$ws.Range("J19").Value = 850
$ws.Range("L19").Value = 850
$ws.Range("N19").Value = -1190

# row 22
$ws.Range("H22").Value = 1161.5555
$ws.Range("I22").Value = 799
$ws.Range("J22").Value = 1265.1428
$ws.Range("K22").Value = 799
$ws.Range("L22").Value = 1265.1428
$ws.Range("M22").Value = -504
$ws.Range("N22").Value = -1855.1428

# row 27
$ws.Range("H27").Value = 1161.5555
$ws.Range("I27").Value = 799
$ws.Range("J27").Value = 1265.1428
$ws.Range("K27").Value = 799
$ws.Range("L27").Value = 1265.1428
$ws.Range("M27").Value = -692
$ws.Range("N27").Value = -1479.1428

# row 40
$ws.Range("H40").Value = 1975.5714
$ws.Range("I40").Value = 1973.3334
$ws.Range("K40").Value = 1973.3334
$ws.Range("M40").Value = -1837.3334

# row 127
$ws.Range("H127").Value = 77497.75
$ws.Range("J127").Value = 77497.75
$ws.Range("L127").Value = 77497.75
$ws.Range("N127").Value = -87417.75

# row 132
$ws.Range("H132").Value = 15000

$ws = $wb.Worksheets.Item("WVR")
# row 107
$ws.Range("H107").Value = 197.33333
$ws.Range("I107").Value = 197.33333
$ws.Range("K107").Value = 591.99999
$ws.Range("M107").Value = 1328.00001

# row 132
$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470
